$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 583, shifting the existing data (rows 583:655) down to 586:658.
$ws.Rows("583:585").Insert()

# Row 583: new weekly record (Zafiro rojo, Primera)
$ws.Cells.Item(583,1).Value = 1
$ws.Cells.Item(583,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(583,3).Value = "Arica y Parinacota"
$ws.Cells.Item(583,4).Value = 44748
$ws.Cells.Item(583,5).Value = 15
$ws.Cells.Item(583,6).Value = 100112002
$ws.Cells.Item(583,7).Value = "Pimiento"
$ws.Cells.Item(583,8).Value = "Zafiro rojo"
$ws.Cells.Item(583,9).Value = "Primera"
$ws.Cells.Item(583,10).Value = 120
$ws.Cells.Item(583,11).Value = 27000
$ws.Cells.Item(583,12).Value = 28000
$ws.Cells.Item(583,13).Value = 27500
$ws.Cells.Item(583,14).Value = "$/caja 15 kilos"
$ws.Cells.Item(583,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(583,16).Value = 1833
$ws.Cells.Item(583,17).Value = 15
$ws.Cells.Item(583,18).Value = "Hortaliza"

# Row 584: new weekly record (Zafiro rojo, Segunda)
$ws.Cells.Item(584,1).Value = 1
$ws.Cells.Item(584,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(584,3).Value = "Arica y Parinacota"
$ws.Cells.Item(584,4).Value = 44748
$ws.Cells.Item(584,5).Value = 15
$ws.Cells.Item(584,6).Value = 100112002
$ws.Cells.Item(584,7).Value = "Pimiento"
$ws.Cells.Item(584,8).Value = "Zafiro rojo"
$ws.Cells.Item(584,9).Value = "Segunda"
$ws.Cells.Item(584,10).Value = 135
$ws.Cells.Item(584,11).Value = 25000
$ws.Cells.Item(584,12).Value = 26000
$ws.Cells.Item(584,13).Value = 25481
$ws.Cells.Item(584,14).Value = "$/caja 15 kilos"
$ws.Cells.Item(584,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(584,16).Value = 1699
$ws.Cells.Item(584,17).Value = 15
$ws.Cells.Item(584,18).Value = "Hortaliza"

# Row 585: new weekly record (Zafiro rojo, Tercera)
$ws.Cells.Item(585,1).Value = 1
$ws.Cells.Item(585,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(585,3).Value = "Arica y Parinacota"
$ws.Cells.Item(585,4).Value = 44748
$ws.Cells.Item(585,5).Value = 15
$ws.Cells.Item(585,6).Value = 100112002
$ws.Cells.Item(585,7).Value = "Pimiento"
$ws.Cells.Item(585,8).Value = "Zafiro rojo"
$ws.Cells.Item(585,9).Value = "Tercera"
$ws.Cells.Item(585,10).Value = 150
$ws.Cells.Item(585,11).Value = 22000
$ws.Cells.Item(585,12).Value = 23000
$ws.Cells.Item(585,13).Value = 22500
$ws.Cells.Item(585,14).Value = "$/caja 15 kilos"
$ws.Cells.Item(585,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(585,16).Value = 1500
$ws.Cells.Item(585,17).Value = 15
$ws.Cells.Item(585,18).Value = "Hortaliza"

Write-Output "Done. UsedRange rows:"
Write-Output $ws.UsedRange.Rows.Count
